# Auto-generated script applying numeric corrections to the Masamune_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(95, 8).Value2 = 37484.8
$ws.Cells.Item(95, 10).Value2 = 37484.8
$ws.Cells.Item(95, 12).Value2 = 37484.8
$ws.Cells.Item(95, 14).Value2 = -42976.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(108, 8).Value2 = 45418
$ws.Cells.Item(108, 10).Value2 = 45418
$ws.Cells.Item(108, 12).Value2 = 45418
$ws.Cells.Item(108, 14).Value2 = -53098

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(109, 8).Value2 = 34626.668
$ws.Cells.Item(109, 10).Value2 = 34626.668
$ws.Cells.Item(109, 12).Value2 = 34626.668
$ws.Cells.Item(109, 14).Value2 = -37400.668

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(126, 8).Value2 = 39497.855
$ws.Cells.Item(126, 10).Value2 = 39497.855
$ws.Cells.Item(126, 12).Value2 = 39497.855
$ws.Cells.Item(126, 14).Value2 = -49377.855

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value2 = 35271.54
$ws.Cells.Item(133, 10).Value2 = 35271.54
$ws.Cells.Item(133, 12).Value2 = 35271.54
$ws.Cells.Item(133, 14).Value2 = -45391.54

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(43, 8).Value2 = 1435612.1
$ws.Cells.Item(43, 9).Value2 = 7600
$ws.Cells.Item(43, 10).Value2 = 1673614.1
$ws.Cells.Item(43, 11).Value2 = 7600
$ws.Cells.Item(43, 12).Value2 = 1673614.1
$ws.Cells.Item(43, 13).Value2 = -7287
$ws.Cells.Item(43, 14).Value2 = -1674240.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value2 = 46704.89
$ws.Cells.Item(80, 10).Value2 = 46704.89
$ws.Cells.Item(80, 12).Value2 = 46704.89
$ws.Cells.Item(80, 14).Value2 = -48700.89

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(83, 8).Value2 = 46704.89
$ws.Cells.Item(83, 10).Value2 = 46704.89
$ws.Cells.Item(83, 12).Value2 = 140114.67
$ws.Cells.Item(83, 14).Value2 = -150098.67

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(95, 8).Value2 = 36935.5
$ws.Cells.Item(95, 10).Value2 = 36935.5
$ws.Cells.Item(95, 12).Value2 = 36935.5
$ws.Cells.Item(95, 14).Value2 = -42427.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(101, 8).Value2 = 48060
$ws.Cells.Item(101, 10).Value2 = 48060
$ws.Cells.Item(101, 12).Value2 = 48060
$ws.Cells.Item(101, 14).Value2 = -54550

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(104, 8).Value2 = 42741.668
$ws.Cells.Item(104, 10).Value2 = 42741.668
$ws.Cells.Item(104, 12).Value2 = 42741.668
$ws.Cells.Item(104, 14).Value2 = -49729.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(105, 8).Value2 = 49362
$ws.Cells.Item(105, 10).Value2 = 49362
$ws.Cells.Item(105, 12).Value2 = 49362
$ws.Cells.Item(105, 14).Value2 = -56350

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(109, 8).Value2 = 39674.6
$ws.Cells.Item(109, 10).Value2 = 39674.6
$ws.Cells.Item(109, 12).Value2 = 39674.6
$ws.Cells.Item(109, 14).Value2 = -42448.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(113, 8).Value2 = 46615
$ws.Cells.Item(113, 10).Value2 = 46615
$ws.Cells.Item(113, 12).Value2 = 46615
$ws.Cells.Item(113, 14).Value2 = -55293

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(114, 8).Value2 = 45961.332
$ws.Cells.Item(114, 10).Value2 = 45961.332
$ws.Cells.Item(114, 12).Value2 = 45961.332
$ws.Cells.Item(114, 14).Value2 = -54639.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(117, 8).Value2 = 37940.3
$ws.Cells.Item(117, 10).Value2 = 37940.3
$ws.Cells.Item(117, 12).Value2 = 37940.3
$ws.Cells.Item(117, 14).Value2 = -47118.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(118, 8).Value2 = 49401
$ws.Cells.Item(118, 10).Value2 = 49401
$ws.Cells.Item(118, 12).Value2 = 49401
$ws.Cells.Item(118, 14).Value2 = -52715

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(125, 8).Value2 = 46282.5
$ws.Cells.Item(125, 10).Value2 = 46282.5
$ws.Cells.Item(125, 12).Value2 = 46282.5
$ws.Cells.Item(125, 14).Value2 = -56122.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(138, 8).Value2 = 44500
$ws.Cells.Item(138, 10).Value2 = 44500
$ws.Cells.Item(138, 12).Value2 = 44500
$ws.Cells.Item(138, 14).Value2 = -54780

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(100, 8).Value2 = 38903.75
$ws.Cells.Item(100, 10).Value2 = 38903.75
$ws.Cells.Item(100, 12).Value2 = 38903.75
$ws.Cells.Item(100, 14).Value2 = -41067.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(112, 8).Value2 = 46661
$ws.Cells.Item(112, 10).Value2 = 46661
$ws.Cells.Item(112, 12).Value2 = 46661
$ws.Cells.Item(112, 14).Value2 = -49615

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(116, 8).Value2 = 45718
$ws.Cells.Item(116, 10).Value2 = 45718
$ws.Cells.Item(116, 12).Value2 = 45718
$ws.Cells.Item(116, 14).Value2 = -54896

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(119, 8).Value2 = 48753
$ws.Cells.Item(119, 10).Value2 = 48753
$ws.Cells.Item(119, 12).Value2 = 48753
$ws.Cells.Item(119, 14).Value2 = -58429

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(126, 8).Value2 = 50772
$ws.Cells.Item(126, 10).Value2 = 50772
$ws.Cells.Item(126, 12).Value2 = 50772
$ws.Cells.Item(126, 14).Value2 = -60652

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(133, 8).Value2 = 48999.75
$ws.Cells.Item(133, 10).Value2 = 48999.75
$ws.Cells.Item(133, 12).Value2 = 48999.75
$ws.Cells.Item(133, 14).Value2 = -59119.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(111, 8).Value2 = 48694
$ws.Cells.Item(111, 10).Value2 = 48694
$ws.Cells.Item(111, 12).Value2 = 48694
$ws.Cells.Item(111, 14).Value2 = -56874

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(116, 8).Value2 = 49370.5
$ws.Cells.Item(116, 10).Value2 = 49370.5
$ws.Cells.Item(116, 12).Value2 = 49370.5
$ws.Cells.Item(116, 14).Value2 = -58548.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(137, 8).Value2 = 47096.668
$ws.Cells.Item(137, 10).Value2 = 47096.668
$ws.Cells.Item(137, 12).Value2 = 47096.668
$ws.Cells.Item(137, 14).Value2 = -57296.668

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(104, 8).Value2 = 46665.25
$ws.Cells.Item(104, 10).Value2 = 46665.25
$ws.Cells.Item(104, 12).Value2 = 46665.25
$ws.Cells.Item(104, 14).Value2 = -53653.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(105, 8).Value2 = 40790.8
$ws.Cells.Item(105, 10).Value2 = 40790.8
$ws.Cells.Item(105, 12).Value2 = 40790.8
$ws.Cells.Item(105, 14).Value2 = -47778.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(110, 8).Value2 = 49459
$ws.Cells.Item(110, 10).Value2 = 49459
$ws.Cells.Item(110, 12).Value2 = 49459
$ws.Cells.Item(110, 14).Value2 = -57639

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(114, 8).Value2 = 38124
$ws.Cells.Item(114, 10).Value2 = 38124
$ws.Cells.Item(114, 12).Value2 = 38124
$ws.Cells.Item(114, 14).Value2 = -46802

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(130, 8).Value2 = 34481.145
$ws.Cells.Item(130, 10).Value2 = 34481.145
$ws.Cells.Item(130, 12).Value2 = 34481.145
$ws.Cells.Item(130, 14).Value2 = -44521.145

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(103, 8).Value2 = 49590
$ws.Cells.Item(103, 10).Value2 = 49590
$ws.Cells.Item(103, 12).Value2 = 49590
$ws.Cells.Item(103, 14).Value2 = -51934

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(110, 8).Value2 = 41137
$ws.Cells.Item(110, 10).Value2 = 41137
$ws.Cells.Item(110, 12).Value2 = 41137
$ws.Cells.Item(110, 14).Value2 = -49317

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(111, 8).Value2 = 45379
$ws.Cells.Item(111, 10).Value2 = 45379
$ws.Cells.Item(111, 12).Value2 = 45379
$ws.Cells.Item(111, 14).Value2 = -53559

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(134, 8).Value2 = 49494.75
$ws.Cells.Item(134, 10).Value2 = 49494.75
$ws.Cells.Item(134, 12).Value2 = 49494.75
$ws.Cells.Item(134, 14).Value2 = -59634.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value2 = 45992
$ws.Cells.Item(16, 10).Value2 = 45992
$ws.Cells.Item(16, 12).Value2 = 45992
$ws.Cells.Item(16, 14).Value2 = -46576

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value2 = 61848.4
$ws.Cells.Item(46, 10).Value2 = 61848.4
$ws.Cells.Item(46, 12).Value2 = 61848.4
$ws.Cells.Item(46, 14).Value2 = -62310.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(105, 8).Value2 = 45664
$ws.Cells.Item(105, 10).Value2 = 45664
$ws.Cells.Item(105, 12).Value2 = 45664
$ws.Cells.Item(105, 14).Value2 = -52652

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(117, 8).Value2 = 40987.6
$ws.Cells.Item(117, 10).Value2 = 40987.6
$ws.Cells.Item(117, 12).Value2 = 40987.6
$ws.Cells.Item(117, 14).Value2 = -50165.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value2 = 48632
$ws.Cells.Item(119, 10).Value2 = 48632
$ws.Cells.Item(119, 12).Value2 = 48632
$ws.Cells.Item(119, 14).Value2 = -58308

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(121, 8).Value2 = 42205.332
$ws.Cells.Item(121, 10).Value2 = 42205.332
$ws.Cells.Item(121, 12).Value2 = 42205.332
$ws.Cells.Item(121, 14).Value2 = -45699.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(127, 8).Value2 = 42429
$ws.Cells.Item(127, 10).Value2 = 42429
$ws.Cells.Item(127, 12).Value2 = 42429
$ws.Cells.Item(127, 14).Value2 = -52349

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(134, 8).Value2 = 61848.4
$ws.Cells.Item(134, 10).Value2 = 61848.4
$ws.Cells.Item(134, 12).Value2 = 185545.2
$ws.Cells.Item(134, 14).Value2 = -190615.2
